$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.702.17"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.694.33"

$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3953"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.75%  "

$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08961"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.274"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.58"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.044"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001325"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").Value = "1.698.44"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.99"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07036"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.22%  "

$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("D24").Value = "24.691.54"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.296"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.361"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.19%  "

$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.185"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.498"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08650"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.054"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.070"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.39"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2740"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.884"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.84%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.47"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09261"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02725"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.474"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7665"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.593"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7171"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.320"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07983"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "
